# Updates the cryptos price/volume snapshot on Sheet1 (and two row swaps:
# Monero/ARBITRUM at rows 47-48, Filecoin->FirstDigitalUSD at row 51) to
# match the refreshed GitHub Actions data pull.
#
# Note: several "Price" values in column D are plain decimal-looking
# strings (e.g. "235.42"); Excel's COM layer auto-converts such strings to
# numbers on assignment. Forcing NumberFormat to "@" (Text) before the
# assignment keeps them as text like the rest of the column, and resetting
# the Style back to "Normal" afterwards avoids leaving a custom number
# format applied to the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '90.105.75'
$ws.Range('E2').Value = '  -1.07%  '
$ws.Range('D3').Value = '3.096.48'
$ws.Range('E3').Value = '  -2.27%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '235.42'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +8.84%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '620.06'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.50%  '
$ws.Range('E7').Value = '  -12.40%  '
$ws.Range('E8').Value = '  -3.39%  '
$ws.Range('E9').Value = '  +0.09%  '
$ws.Range('D10').Value = '3.093.64'
$ws.Range('E10').Value = '  -2.28%  '
$ws.Range('E11').Value = '  -5.70%  '
$ws.Range('E12').Value = '  -3.04%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000250'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.08%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '35.32'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.38%  '
$ws.Range('D15').Value = '89.858.62'
$ws.Range('E15').Value = '  -1.07%  '
$ws.Range('E16').Value = '  -7.04%  '
$ws.Range('D17').Value = '3.668.67'
$ws.Range('E17').Value = '  -2.25%  '
$ws.Range('D18').Value = '3.077.87'
$ws.Range('E18').Value = '  -3.04%  '
$ws.Range('E19').Value = '  +1.45%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0000213'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.38%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.85'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -6.00%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '433.32'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -10.01%  '
$ws.Range('E23').Value = '  +5.71%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '8.82'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -3.70%  '
$ws.Range('E25').Value = '  -2.99%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '86.35'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -11.16%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.83'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -4.38%  '
$ws.Range('D28').Value = '3.274.11'
$ws.Range('E28').Value = '  -1.97%  '
$ws.Range('E29').Value = '  +0.06%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.11'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.34%  '
$ws.Range('E31').Value = '  +0.08%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.158'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.34%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.192'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.85%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '25.65'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -9.13%  '
$ws.Range('E35').Value = '  +3.84%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.73'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.52%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '7.15'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.27%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '497.51'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -5.43%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.88'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.15%  '
$ws.Range('E40').Value = '  -3.22%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.66'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +58.42%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0869'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -6.19%  '
$ws.Range('E43').Value = '  -0.58%  '
$ws.Range('E44').Value = '  -0.01%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.398'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -5.90%  '
$ws.Range('E46').Value = '  -6.28%  '
$ws.Range('B47').Value = 'ARBITRUM'
$ws.Range('C47').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.682'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -4.89%  '
$ws.Range('B48').Value = 'Monero'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '151.34'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.59%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '44.44'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.49%  '
$ws.Range('E50').Value = '  -4.49%  '
$ws.Range('B51').Value = 'FirstDigitalUSD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.00'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.11%  '
